$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from SCD0183 to SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID cell (B2) from "DGS-198" to "SCD0011-014"
$ws.Range("B2").Value = "SCD0011-014"

# Move/restore the active selection to B3, matching the saved view state.
$ws.Range("B3").Select()
